$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-09-15"

# Update the September label (row 10, column A) to new "through" date
$ws.Range("A10").Value = "September (through 09-15)"

# July (row 8) - update 2022 column (I)
$ws.Range("I8").Value = 164

# September (row 10) - update 2016..2021 columns (C..H); I10 (2022) also changes
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = 37
$ws.Range("E10").Value = 29
$ws.Range("F10").Value = 32
$ws.Range("G10").Value = 57
$ws.Range("H10").Value = 78
$ws.Range("I10").Value = 67

# Total (row 11) - update 2016..2021 columns (C..H); I11 (2022) unchanged
$ws.Range("C11").Value = 407
$ws.Range("D11").Value = 588
$ws.Range("E11").Value = 519
$ws.Range("F11").Value = 381
$ws.Range("G11").Value = 841
$ws.Range("H11").Value = 1148

$wb.Save()
